$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.4
$ws.Range("I2").Value = 6.25
$ws.Range("M2").Value = 1.02
$ws.Range("N2").Value = 11
$ws.Range("Z2").Value = 9
$ws.Range("AC2").Value = 11
$ws.Range("AK2").Value = 81
$ws.Range("AQ2").Value = 19
$ws.Range("AW2").Value = 8.5
$ws.Range("AX2").Value = 41
$ws.Range("AZ2").Value = 151
